$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.539
$ws.Range("C14").Value = -11.944
$ws.Range("C21").Value = -12.953
$ws.Range("C23").Value = -13.262
$ws.Range("C25").Value = -12.763
